# edit.ps1 - PowerPoint COM-interop script (iron_native / run_com)
#
# This reproduces two logical changes from the source commit:
#
#   1. Three tables (on the slides holding the "Google Shape;202/…;p26"
#      table graphic frames) get their <a:tableStyleId> switched from the
#      custom "Table_0" style ({0D7A0236-2E5D-44D4-8565-69FD1E0A0D33}) to
#      the built-in "No Style, No Grid" style
#      ({F8D37DC5-4563-4A5A-8862-D0D3242EEEB7}).
#
#   2. The deck's two DrawingML themes (ppt/theme/theme1.xml used by the
#      slide master, and ppt/theme/theme2.xml used by the notes master)
#      get their color palettes swapped: theme1 (currently the "Integral"
#      / "Red Violet" palette) becomes the "Office" palette, and theme2
#      (currently "Office") becomes "Integral" / "Red Violet".
#
# Helper: build the COM "RGB" style integer (0xBBGGRR packing used by the
# PowerPoint object model's ColorFormat/ColorScheme .RGB property) from a
# "RRGGBB" hex string.
function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style swap on the three affected slides.
# ---------------------------------------------------------------------
$newTableStyleId = "{F8D37DC5-4563-4A5A-8862-D0D3242EEEB7}"

for ($slideIdx = 14; $slideIdx -le 16; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shpIdx = 1; $shpIdx -le $slide.Shapes.Count; $shpIdx++) {
        $shape = $slide.Shapes.Item($shpIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme color-scheme swap.
#
# The host's Table/Master/Slide ".ColorScheme" accessor always edits the
# single DrawingML theme that backs the slide master (ppt/theme/theme1.xml)
# regardless of which object it is read from, so it is used here to push
# theme1's 12 scheme colors (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) to what used to be theme2's ("Office") palette.
# ---------------------------------------------------------------------
$officePalette = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$scheme = $p.SlideMaster.ColorScheme
for ($i = 0; $i -lt $officePalette.Length; $i++) {
    $scheme.Colors($i + 1).RGB = HexToComRgb $officePalette[$i]
}
